$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 252; this pushes the existing rows 252-256
# (old data) down to rows 255-259, matching the target layout.
$ws.Rows.Item(252).Insert()
$ws.Rows.Item(252).Insert()
$ws.Rows.Item(252).Insert()

# Copy the date cell style (numFmt) from the row above into the new rows'
# D column so the new date values keep the same date formatting.
$ws.Range("D251").Copy()
$ws.Range("D252:D254").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 252: Cuatro cascos verde
$ws.Range("A252").Value = 11
$ws.Range("B252").Value = "Vega Monumental Concepción"
$ws.Range("C252").Value = "Bíobío"
$ws.Range("D252").Value = 44628
$ws.Range("E252").Value = 8
$ws.Range("F252").Value = 100112002
$ws.Range("G252").Value = "Pimiento"
$ws.Range("H252").Value = "Cuatro cascos verde"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 180
$ws.Range("K252").Value = 9500
$ws.Range("L252").Value = 10000
$ws.Range("M252").Value = 9722
$ws.Range("N252").Value = "$/caja 15 kilos"
$ws.Range("O252").Value = "Provincia de Limarí"
$ws.Range("P252").Value = 648
$ws.Range("Q252").Value = 15
$ws.Range("R252").Value = "Hortaliza"

# Row 253: Morrón rojo
$ws.Range("A253").Value = 11
$ws.Range("B253").Value = "Vega Monumental Concepción"
$ws.Range("C253").Value = "Bíobío"
$ws.Range("D253").Value = 44628
$ws.Range("E253").Value = 8
$ws.Range("F253").Value = 100112002
$ws.Range("G253").Value = "Pimiento"
$ws.Range("H253").Value = "Morrón rojo"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 220
$ws.Range("K253").Value = 15000
$ws.Range("L253").Value = 16000
$ws.Range("M253").Value = 15545
$ws.Range("N253").Value = "$/caja 15 kilos"
$ws.Range("O253").Value = "Provincia de Limarí"
$ws.Range("P253").Value = 1036
$ws.Range("Q253").Value = 15
$ws.Range("R253").Value = "Hortaliza"

# Row 254: Zafiro rojo
$ws.Range("A254").Value = 11
$ws.Range("B254").Value = "Vega Monumental Concepción"
$ws.Range("C254").Value = "Bíobío"
$ws.Range("D254").Value = 44628
$ws.Range("E254").Value = 8
$ws.Range("F254").Value = 100112002
$ws.Range("G254").Value = "Pimiento"
$ws.Range("H254").Value = "Zafiro rojo"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 180
$ws.Range("K254").Value = 25000
$ws.Range("L254").Value = 26000
$ws.Range("M254").Value = 25444
$ws.Range("N254").Value = "$/caja 15 kilos"
$ws.Range("O254").Value = "Región de Arica y Parinacota"
$ws.Range("P254").Value = 1696
$ws.Range("Q254").Value = 15
$ws.Range("R254").Value = "Hortaliza"
